$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-42 (date shift +3 days, new flow values) and add new rows 43-48
$data = @(
    @(2, 45828, 1.741, 1.894),
    @(3, 45828.01041666666, 20.577, 0),
    @(4, 45828.02083333334, 19.888, 0.006),
    @(5, 45828.03125, 11.107, 0.005),
    @(6, 45828.04166666666, 20.861, 0),
    @(7, 45828.05208333334, 2.288, 14.984),
    @(8, 45828.0625, 0, 26.582),
    @(9, 45828.07291666666, 1.578, 4.71),
    @(10, 45828.08333333334, 0.375, 20.245),
    @(11, 45828.09375, 0, 32.981),
    @(12, 45828.10416666666, 0, 88.02),
    @(13, 45828.11458333334, 0, 50.052),
    @(14, 45828.125, 0, 9.516999999999999),
    @(15, 45828.13541666666, 0, 12.381),
    @(16, 45828.14583333334, 0, 22.014),
    @(17, 45828.15625, 1.933, 12.755),
    @(18, 45828.16666666666, 32.357, 0),
    @(19, 45828.17708333334, 2.4, 9.449),
    @(20, 45828.1875, 0, 43.29),
    @(21, 45828.19791666666, 12.418, 4.165),
    @(22, 45828.20833333334, 37.803, 0),
    @(23, 45828.21875, 12.373, 1.644),
    @(24, 45828.22916666666, 4.815, 3.338),
    @(25, 45828.23958333334, 8.316000000000001, 0.372),
    @(26, 45828.25, 16.048, 0),
    @(27, 45828.26041666666, 26.585, 0),
    @(28, 45828.27083333334, 0.229, 5.683),
    @(29, 45828.28125, 0, 10.126),
    @(30, 45828.29166666666, 11.719, 0.617),
    @(31, 45828.30208333334, 21.521, 0.037),
    @(32, 45828.3125, 1.416, 5.532),
    @(33, 45828.32291666666, 2.158, 2.7),
    @(34, 45828.33333333334, 2.332, 8.907),
    @(35, 45828.34375, 0, 7.674),
    @(36, 45828.35416666666, 0.007, 9.305),
    @(37, 45828.36458333334, 0, 8.778),
    @(38, 45828.375, 0, 49.621),
    @(39, 45828.38541666666, 0, 60.767),
    @(40, 45828.39583333334, 0, 26.779),
    @(41, 45828.40625, 0, 2.016),
    @(42, 45828.41666666666, 0.015, 14.557),
    @(43, 45828.42708333334, 0, 19.122),
    @(44, 45828.4375, 1.469, 1.726),
    @(45, 45828.44791666666, 8.16, 3.712),
    @(46, 45828.45833333334, 0, 11.97),
    @(47, 45828.46875, 0, 7.164),
    @(48, 45828.47916666666, 6.163, 0.707)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
}

# Ensure date/time number format (style) applies to the newly added column-A cells (rows 43-48)
$ws.Range("A43:A48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
